# Add a second row of "test data" below the existing header row.
# The new values are written with a leading apostrophe so Excel stores
# them as literal text ("5000.0" / "6000.0") instead of re-interpreting
# them as numbers (which would drop the trailing ".0").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "'5000.0"
$ws.Range("B2").Value = "'6000.0"
